$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 = repaymentstrategy field. Change its value from "RBI (India)" to
# the new repayment strategy scenario description.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Update the active selection to reflect where the edit was made.
$ws.Range("B17").Select()
